# angelou.docx — "Changed the main contents as part 9"
#
# 1) Paragraph 1: keep the original sentence, add two trailing spaces, then
#    append a red parenthetical remark split across three runs (matching how
#    it was evidently typed/edited in separate passes):
#       "(This is a change – Ve" / "rsion for main branch" / ")"
# 2) Paragraph 3 (the empty Menlo-formatted paragraph) becomes a bare, fully
#    empty paragraph with no paragraph/run formatting at all.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: text + colored parenthetical
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
# Range that covers the paragraph's content but excludes the trailing
# paragraph mark, so InsertAfter keeps landing inside paragraph 1.
$body1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)

# Two trailing spaces after the existing sentence (still default/black).
$body1.InsertAfter("  ")

# "(This is a change – Ve" in red.
$chunkStart = $body1.End
$body1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run1 = $d.Range($chunkStart, $body1.End)
$run1.Font.Color = 255

# "rsion for main branch" in red (separate run).
$chunkStart = $body1.End
$body1.InsertAfter("rsion for main branch")
$run2 = $d.Range($chunkStart, $body1.End)
$run2.Font.Color = 255

# ")" in red (separate run).
$chunkStart = $body1.End
$body1.InsertAfter(")")
$run3 = $d.Range($chunkStart, $body1.End)
$run3.Font.Color = 255

# ---------------------------------------------------------------------
# 2) Third paragraph: strip it down to a totally bare empty paragraph
# ---------------------------------------------------------------------
# Insert a brand-new, unformatted paragraph right after paragraph 2 (the
# plain "It will be treated..." paragraph), so it doesn't inherit any of
# paragraph 3's Menlo/baseline formatting.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()

# The freshly inserted paragraph (now paragraph 3) still carries an
# internal empty run recording "current" insertion formatting. Typing a
# character into it and deleting it again collapses that run away,
# leaving a truly empty <w:p/>.
$p3 = $d.Paragraphs(3)
$p3.Range.InsertAfter("X")
$stray = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$stray.Delete()

# Finally, remove the old Menlo-formatted empty paragraph, which has been
# pushed down to index 4 by the insertion above.
$d.Paragraphs(4).Range.Delete()
